# Token/Exit System Settings Changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Place, Phone No) before the "Material" column (old column G).
$ws.Columns("G:H").Insert()

# Insert one new column (Round off) before the "Nett Wt" column
# (old column O, now at column Q after the previous insert).
$ws.Columns("Q:Q").Insert()

# Set header text for the newly inserted columns.
$ws.Range("G1").Value = "Place"
$ws.Range("H1").Value = "Phone No"
$ws.Range("Q1").Value = "Round off"

# Give the new columns an explicit, reasonable best-fit-style width (matches
# the narrower layout the workbook now uses for its header labels).
$ws.Columns("G").ColumnWidth = 4.498697916666667
$ws.Columns("H").ColumnWidth = 8.166666666666666
$ws.Columns("Q").ColumnWidth = 8.276041666666666

# Clear the stray "0" placeholder values that are no longer populated in
# the data rows for "No Of Bags", "Charges", "Bag Deduction", "Final Wt"
# and "Final Amount" columns (shifted to J, K, P, T, U respectively).
$ws.Range("J2:K4").ClearContents()
$ws.Range("P2:P4").ClearContents()
$ws.Range("T2:U4").ClearContents()

# Restore the selection Excel leaves behind after this kind of edit.
$ws.Range("A5:XFD17").Select()
